$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Test Case No.12 (row 13) 즐겨찾기 now marked as Pass ("O") in the
# 통과 여부 (Pass/Fail) column G, matching the formatting already used
# by the other cells in that column (center aligned, wrapped text).
$g13 = $ws.Range("G13")
$g13.HorizontalAlignment = -4108
$g13.VerticalAlignment = -4108
$g13.WrapText = $true
$g13.Value = "O"

# Reflect the updated selection / scroll position recorded in the saved view state.
$ws.Range("A11").Select()
$ws.Range("F14").Select()
